$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised historical values ---
$ws.Range("C289:F289").Value = 32957600000
$ws.Range("C322:F322").Value = 36383400000
$ws.Range("C327:F327").Value = 37472200000
$ws.Range("C328:F328").Value = 37492700000

# --- Append two new data rows (352, 353), matching formatting of row 351 ---
$ws.Range("A351").Copy()
$ws.Range("A352").PasteSpecial(-4122)
$ws.Range("A353").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A352").Value = 44986.45833333334
$ws.Range("B352").Value = "ECONOMICS:JOM2"
$ws.Range("C352:F352").Value = 41710700000
$ws.Range("G352").Value = 0

$ws.Range("A353").Value = 45017.45833333334
$ws.Range("B353").Value = "ECONOMICS:JOM2"
$ws.Range("C353:F353").Value = 41839100000
$ws.Range("G353").Value = 0
